# Moved some configuration settings to config and added support for
# multiple sequence passes.
#
# - NormalSearch (sheet2): drop the trailing "CheckDB" row.
# - FeelingLuckySearch (sheet3): simplify each step row to match the
#   NormalSearch pattern (target/value now sourced from config, literal
#   numeric "Value" for Wait steps), dropping the old Target/Value/
#   Description trio for most rows.
# - Both sheets' active-cell selection moves to C11 (room for the extra
#   pass rows).

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("NormalSearch")
$ws3 = $wb.Worksheets.Item("FeelingLuckySearch")

# --- NormalSearch: remove the CheckDB row (row 6) ---
$ws2.Rows.Item(6).Delete()

# --- FeelingLuckySearch: rework rows 2-5 ---
# Row 2: Type -> map(searchtextbox2) / data(searchkeyword), drop Description
$ws3.Range("C2").Value = "`$map(searchpage|searchtextbox2)"
$ws3.Range("D2").Value = "`$data(searchkeyword)"
$ws3.Range("E2").ClearContents()

# Row 3: now a Wait step with a literal pass-count value
$ws3.Range("A3").Value = "y"
$ws3.Range("B3").Value = "Wait"
$ws3.Range("C3").ClearContents()
$ws3.Range("D3").Value = 1
$ws3.Range("E3").ClearContents()

# Row 4: Click -> map(searchbutton)
$ws3.Range("A4").Value = "y"
$ws3.Range("B4").Value = "Click"
$ws3.Range("C4").Value = "`$map(searchpage|searchbutton)"
$ws3.Range("D4").ClearContents()
$ws3.Range("E4").ClearContents()

# Row 5: another Wait step, second pass count
$ws3.Range("A5").Value = "y"
$ws3.Range("B5").Value = "Wait"
$ws3.Range("C5").ClearContents()
$ws3.Range("D5").Value = 2
$ws3.Range("E5").ClearContents()

# --- Selection bookkeeping: both sheets now select C11 ---
# Select FeelingLuckySearch's cell first so the workbook's "active" tab
# ends up back on NormalSearch (matching the original tabSelected sheet).
$ws3.Range("C11").Select()
$ws2.Range("C11").Select()
